# New crime data collected — weekly refresh of the 68th Precinct CompStat sheet.
# Updates: report header (volume/number + date range), and the weekly/28-day/
# YTD complaint figures + percent-change columns for rows 16-28 and 33, some of
# which flip between a numeric count and the sheet's "0" / "***.*" text
# placeholders (used when a count or percent-change is not meaningful).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text (shared-string runs collapsed to their plain concatenation;
# all runs in each string share identical formatting, so this is a faithful
# edit of "Volume 32   Number  23/24" and the "Report Covering the Week ..."
# line).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/9/2025  Through  6/15/2025"

# ---------------------------------------------------------------------------
# Helper donor cells whose style we reuse via Copy + PasteSpecial(xlPasteFormats)
# so that cells which change between a numeric style and the plain/text style
# pick up the exact same cellXf the rest of the sheet already uses (instead of
# Excel minting a brand-new style). None of these donor cells are themselves
# edited below, and PasteSpecial(formats) never touches the destination value.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122
$styleTextDonor = "C23"   # s=13 (General / plain text style)
$styleCountDonor = "I19"  # s=14 (#,##0 style used for counts)
$stylePctDonor = "K19"    # s=15 (#,##0.0 style used for percents)

function Set-Count {
    param($ref, $value)
    $ws.Range($ref).Value = $value
}

function Set-Pct {
    param($ref, $value)
    $ws.Range($ref).Value = $value
}

function Set-TextPlaceholder {
    param($ref, $text)
    # Force the literal string (rather than letting "0" parse back to a
    # number) by entering it with a leading apostrophe, then immediately
    # reapply the sheet's normal "text" cell format over it so the
    # quote-prefix flag doesn't stick and the style matches the rest of the
    # column (s=13).
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($styleTextDonor).Copy()
    $ws.Range($ref).PasteSpecial($xlPasteFormats)
}

function Set-NumericFromText {
    param($ref, $value, $donor)
    # Currently holds one of the text placeholders ("0" / "***.*"); restyle to
    # the numeric cellXf first, then write the number.
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial($xlPasteFormats)
    $ws.Range($ref).Value = $value
}

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-Count "C16" 1
Set-Pct   "E16" -50
Set-Count "F16" 7
Set-Count "G16" 4
Set-Pct   "H16" 75
Set-Count "I16" 30
Set-Count "J16" 33
Set-Pct   "K16" -9.090909090909
Set-Pct   "L16" -18.918918918918
Set-Pct   "M16" -38.775510204081
Set-Pct   "N16" -89.090909090909

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
Set-Count "C17" 3
Set-Count "D17" 3
Set-Pct   "E17" 0
Set-Count "F17" 12
Set-Count "G17" 13
Set-Pct   "H17" -7.692307692307
Set-Count "I17" 75
Set-Count "J17" 54
Set-Pct   "K17" 38.888888888888
Set-Pct   "L17" 13.636363636363
Set-Pct   "M17" 59.574468085106
Set-Pct   "N17" -39.516129032258

# ---------------------------------------------------------------------------
# Row 18 - Burglary (WTD 2025 count drops to 0 -> rendered as text "0")
# ---------------------------------------------------------------------------
Set-TextPlaceholder "C18" "0"
Set-Count "D18" 3
Set-Pct   "E18" -100
Set-Count "F18" 2
Set-Count "G18" 6
Set-Pct   "H18" -66.666666666666
Set-Count "J18" 45
Set-Pct   "K18" -11.111111111111
Set-Pct   "L18" -27.272727272727
Set-Pct   "M18" -62.962962962963
Set-Pct   "N18" -91.836734693877

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
Set-Count "C19" 5
Set-Count "D19" 6
Set-Pct   "E19" -16.666666666666
Set-Count "F19" 30
Set-Count "G19" 28
Set-Pct   "H19" 7.142857142857
Set-Count "I19" 150
Set-Count "J19" 194
Set-Pct   "K19" -22.680412371134
Set-Pct   "L19" -24.623115577889
Set-Pct   "M19" 11.111111111111
Set-Pct   "N19" -28.229665071770

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-Count "C20" 5
Set-Count "D20" 3
Set-Pct   "E20" 66.666666666666
Set-Count "F20" 8
Set-Pct   "H20" -33.333333333333
Set-Count "I20" 51
Set-Count "J20" 78
Set-Pct   "K20" -34.615384615384
Set-Pct   "L20" 4.081632653061
Set-Pct   "M20" -28.169014084507
Set-Pct   "N20" -94.420131291028

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold row, styles 17/18 - values only)
# ---------------------------------------------------------------------------
Set-Count "C21" 14
Set-Count "D21" 17
Set-Pct   "E21" -17.647058823529
Set-Count "F21" 59
Set-Count "G21" 63
Set-Pct   "H21" -6.349206349206
Set-Count "I21" 353
Set-Count "J21" 407
Set-Pct   "K21" -13.267813267813
Set-Pct   "L21" -14.734299516908
Set-Pct   "M21" -15.347721822542
Set-Pct   "N21" -82.576505429417

# ---------------------------------------------------------------------------
# Row 22 - Transit (WTD counts + %chg all become the "0" / "***.*" text
# placeholders this week)
# ---------------------------------------------------------------------------
Set-TextPlaceholder "C22" "0"
Set-TextPlaceholder "D22" "0"
Set-TextPlaceholder "E22" "***.*"

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
Set-Count "C24" 14
Set-Count "D24" 26
Set-Pct   "E24" -46.153846153846
Set-Count "F24" 55
Set-Count "G24" 102
Set-Pct   "H24" -46.078431372549
Set-Count "I24" 494
Set-Count "J24" 624
Set-Pct   "K24" -20.833333333333
Set-Pct   "L24" -29.428571428571
Set-Pct   "M24" -6.439393939393

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
Set-Count "C25" 4
Set-Count "D25" 16
Set-Pct   "E25" -75
Set-Count "F25" 26
Set-Count "G25" 61
Set-Pct   "H25" -57.377049180327
Set-Count "I25" 242
Set-Count "J25" 398
Set-Pct   "K25" -39.195979899497
Set-Pct   "L25" -33.698630136986

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
Set-Count "C26" 9
Set-Count "D26" 7
Set-Pct   "E26" 28.571428571428
Set-Count "F26" 44
Set-Count "G26" 43
Set-Pct   "H26" 2.325581395348
Set-Count "I26" 180
Set-Count "J26" 166
Set-Pct   "K26" 8.433734939759
Set-Pct   "L26" 2.272727272727
Set-Pct   "M26" 11.111111111111

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape* (WTD 2024 drops to 0 -> text placeholders)
# ---------------------------------------------------------------------------
Set-TextPlaceholder "D27" "0"
Set-TextPlaceholder "E27" "***.*"

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-Count "D28" 2
Set-Count "F28" 1
Set-Count "G28" 4
Set-Pct   "H28" -75
Set-Count "J28" 25
Set-Pct   "K28" -24
Set-Pct   "L28" 72.727272727272

# ---------------------------------------------------------------------------
# Row 33 - Hate Crimes (28-day counts come back to life from the "0" / "***.*"
# placeholders, now populated with real numbers)
# ---------------------------------------------------------------------------
Set-NumericFromText "D33" 2    $styleCountDonor
Set-NumericFromText "E33" -100 $stylePctDonor
Set-NumericFromText "G33" 2    $styleCountDonor
Set-NumericFromText "H33" -100 $stylePctDonor
Set-Count "J33" 4
